$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function New-Run($text, $bold, $italic) {
    return @{ text = $text; bold = $bold; italic = $italic }
}

# Inserts a new paragraph right after the paragraph with index $paraIndex.
# $runsSpec is an array of hashtables produced by New-Run (possibly empty for
# a blank paragraph). $style is a Word style name (or $null for "Normal").
# $leftIndent is a left-indent value expressed in points (or $null/0).
# Returns the index of the newly created paragraph.
function Insert-ParaAfterIndex($paraIndex, $runsSpec, $style, $leftIndent) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $newIndex = $paraIndex + 1
    $newPara = $d.Paragraphs($newIndex)

    if ($style) {
        $newPara.Style = $style
    } else {
        $newPara.Style = "Normal"
    }
    if ($leftIndent) {
        $newPara.LeftIndent = $leftIndent
    }

    if ($runsSpec -and $runsSpec.Count -gt 0) {
        $fullText = ""
        foreach ($r in $runsSpec) {
            $fullText += $r.text
        }
        $newPara.Range.Text = $fullText
        $pos = $newPara.Range.Start
        foreach ($r in $runsSpec) {
            $len = $r.text.Length
            $subRng = $d.Range($pos, $pos + $len)
            $wantBold = 0
            if ($r.bold) { $wantBold = 1 }
            if ($subRng.Font.Bold -ne $wantBold) { $subRng.Font.Bold = $wantBold }
            $wantItalic = 0
            if ($r.italic) { $wantItalic = 1 }
            if ($subRng.Font.Italic -ne $wantItalic) { $subRng.Font.Italic = $wantItalic }
            $pos += $len
        }
    } else {
        if ($newPara.Range.Font.Bold -ne 0) { $newPara.Range.Font.Bold = 0 }
        if ($newPara.Range.Font.Italic -ne 0) { $newPara.Range.Font.Italic = 0 }
    }

    return $newIndex
}

# ---------------------------------------------------------------------------
# 1) "FSC och PEFC" -> "FSC, Chain of Custody, Controlled Wood och PEFC"
# ---------------------------------------------------------------------------

$old1 = "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC och PEFC. I BILAGA 1 finns artfakta om fridlysta arter."
$new1 = "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) New "1.3.1" paragraph right after PRINCIP 1 paragraph
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("Certifikatsinnehavaren ska följa alla tillämpliga lagar, förordningar och nationellt ratificerade internationella avtal, konventioner och överenskommelser.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index

$idx = Insert-ParaAfterIndex $idx @( (New-Run "1.3.1 " $true $false), (New-Run "Tillämpliga lagar och föreskrifter för brukandet av skogen följs." $false $false) ) $null $null

# ---------------------------------------------------------------------------
# 3) Three new paragraphs after the "6.4 Certifikatsinnehavaren..." paragraph
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index

$idx = Insert-ParaAfterIndex $idx @( (New-Run "6.4.1 " $true $false), (New-Run "Följande biotoper undantas från alla skogsbruksåtgärder, förutom åtgärder påkallade för att bevara eller främja biotopens naturliga eller hävdbetingade biologiska mångfald:" $false $false) ) $null $null
$idx = Insert-ParaAfterIndex $idx @( (New-Run "b) nyckelbiotoper enligt Skogsstyrelsens definition och metod (1995)" $false $false) ) $null 17.85
$idx = Insert-ParaAfterIndex $idx @( (New-Run "6.4.1 " $true $false), (New-Run "Bevarandeåtgärder genomförs för de kända förekomster av rödlistade arter som påverkas av skogsbruk." $false $false) ) $null $null

# ---------------------------------------------------------------------------
# 4) Update "I det avverkningsanmälda skogsområdet har..." Kommentar text
# ---------------------------------------------------------------------------

$old4 = "I det avverkningsanmälda skogsområdet har naturvårdsarter och rödlistade arter sina livsmiljöer och växtplatser."
$new4 = "I det avverkningsanmälda skogsområdet har 4 naturvårdsarter varav 4 rödlistade arter sina livsmiljöer och växtplatser."
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Large new block: Chain of Custody / FSC policy / Controlled Wood
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute($new4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Ur Chain of Custody Certification (FSC-STD-40-004 ver 3.0)" $false $false) ) "Heading 2" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "1.3 " $true $false), (New-Run "The organization shall commit to the FSC values as defined in FSC-POL-01-004 Policy for the Association of Organizations with FSC. " $false $false) ) $null $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Ur FSC:s policy för associerade organisationer (FSC-POL-01-004)" $false $false) ) "Heading 2" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Som ”Chain of Custody”-certifierad organisation är skogsbolaget bunden av de fastställda reglerna i Del 1 Punkt 1 c) i Policy för organisationer associerade med FSC (FSC-POL-01-004 V2-0 EN + SVE version 2012-03-02): " $false $false) ) $null $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "De organisationer FSC kan acceptera association med får inte direkt eller indirekt ha några kopplingar till nedanstående, oacceptabla aktiviteter:" $false $false) ) "List Number" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "a) Olaglig skogsavverkning och handel med olagligt avverkat virke eller skogsprodukter" $false $false) ) $null 17.85
$idx = Insert-ParaAfterIndex $idx @( (New-Run "..." $false $false) ) $null 17.85
$idx = Insert-ParaAfterIndex $idx @( (New-Run "c) Skogsbruk som förstör höga naturvärden" $false $false) ) $null 17.85
$idx = Insert-ParaAfterIndex $idx @( (New-Run "d) Betydande omvandling av skog till plantager eller annan, icke skoglig, markanvändning" $false $false) ) $null 17.85

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Kommentar: " $true $true), (New-Run "Avverkning av skog med höga naturvärden samt skada på fridlysta arter strider både mot FSC Controlled Wood-standarden och FSC:s policy for associerade organisationer." $false $true) ) "List Bullet" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Ur FSC Controlled Wood (FSC-STD-40-005)" $false $false) ) "Heading 2" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Virke som inte accepteras i FSC-märkta produkter (oacceptabelt ursprung) enligt FSC Controlled Wood (FSC-STD-40-005):" $false $false) ) $null $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Illegalt avverkat virke." $false $false) ) "List Number" $null
$idx = Insert-ParaAfterIndex $idx @( (New-Run "..." $false $false) ) "List Number" $null
$idx = Insert-ParaAfterIndex $idx @( (New-Run "Virke från avverkningar som hotar höga naturvärden." $false $false) ) "List Number" $null
$idx = Insert-ParaAfterIndex $idx @( (New-Run "Virke från skog som konverteras till plantager eller icke-skogligt bruk." $false $false) ) "List Number" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Kommentar: " $true $true), (New-Run "Avverkning av skog med höga naturvärden samt skada på fridlysta arter strider både mot FSC Controlled Wood-standarden och FSC:s policy for associerade organisationer." $false $true) ) "List Bullet" $null

# ---------------------------------------------------------------------------
# 6) Two blank paragraphs + new closing paragraph after 2nd PEFC Kommentar
# ---------------------------------------------------------------------------

$rngHeading = $d.Content
$rngHeading.Find.Execute("Ur PEFC-standarden gällande lagefterlevnad", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2 = $d.Range($rngHeading.End, $d.Content.End)
$rng2.Find.Execute("I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng2.Paragraphs(1).Index

$idx = Insert-ParaAfterIndex $idx @() $null $null
$idx = Insert-ParaAfterIndex $idx @() $null $null
$idx = Insert-ParaAfterIndex $idx @( (New-Run "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden." $false $false) ) $null $null

# ---------------------------------------------------------------------------
# 7) New "Spillkråka" section after the knärot references, before sectPr
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("SLU Artdatabanken, Uppsala", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Spillkråka – ekologi samt krav på livsmiljön" $false $false) ) "Heading 1" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Spillkråka (NT) är rödlistad som nära hotad, fridlyst enligt §4 Artskyddsförordningen och ingår i bilaga 1 i EU:s fågeldirektiv. Spillkråka lever i både barr- och blandskog liksom i ren lövskog. De tätaste populationerna tenderar att finnas i äldre, variationsrik blandskog med gott om död ved och gamla träd. Varje par utnyttjar 400-1 000 hektar skog beroende på skogens kvalitet. En minskning av populationen pågår på grund av minskad tillgång på lämpliga bo- och födoträd och minskad födotillgång. Spillkråkans minskningstakt har uppgått till 19 (24-10) % under de senaste 15 åren. Skogsbruk med korta omloppstider och täta, homogena ungskogar utgör det största hotet (Artdatabanken 2023)." $false $false) ) $null $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "Referenser - spillkråka" $false $false) ) "Heading 2" $null

$idx = Insert-ParaAfterIndex $idx @( (New-Run "SLU Artdatabanken, 2021. " $false $false), (New-Run "Artfaktablad. Naturvård – artfakta. " $false $true), (New-Run "SLU Artdatabanken, Uppsala" $false $false) ) $null $null

# ---------------------------------------------------------------------------
# 8) Header updates (date + Kopia line)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("2023-10-13", $false, $false, $false, $false, $false, $true, 1, $false, "2023-10-22", 2) | Out-Null
$d.Content.Find.Execute("Kopia: DNV och FSC", $false, $false, $false, $false, $false, $true, 1, $false, "Kopia: Revisor xx och FSC", 2) | Out-Null

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        $hf.Range.Find.Execute("2023-10-13", $false, $false, $false, $false, $false, $true, 1, $false, "2023-10-22", 2) | Out-Null
        $hf.Range.Find.Execute("Kopia: DNV och FSC", $false, $false, $false, $false, $false, $true, 1, $false, "Kopia: Revisor xx och FSC", 2) | Out-Null
    }
}

Write-Output "All edits applied."
